$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A = Oyuncu Adı (player name), Column B = Pozisyon (position),
# Column C = Takım (team), for rows 2..19 (header is row 1).

$data = @(
    @{ Row = 2;  Name = "Andrew Nembhard";          Pos = "PG,SG";      Team = "Indiana Pacers" }
    @{ Row = 3;  Name = "Jordan Poole";              Pos = "PG,SG";      Team = "Washington Wizards" }
    @{ Row = 4;  Name = "Shai Gilgeous-Alexander";   Pos = "PG,SG";      Team = "Oklahoma City Thunder" }
    @{ Row = 5;  Name = "Toumani Camara";            Pos = "SF,PF";      Team = "Portland Trail Blazers" }
    @{ Row = 6;  Name = "Tobias Harris";             Pos = "SF,PF";      Team = "Detroit Pistons" }
    @{ Row = 7;  Name = "Jaylin Williams";           Pos = "PF,C";       Team = "Oklahoma City Thunder" }
    @{ Row = 8;  Name = "Jalen Williams";            Pos = "SG,SF,PF,C"; Team = "Oklahoma City Thunder" }
    @{ Row = 9;  Name = "Cason Wallace";             Pos = "PG,SG";      Team = "Oklahoma City Thunder" }
    @{ Row = 10; Name = "Isaiah Stewart";            Pos = "PF,C";       Team = "Detroit Pistons" }
    @{ Row = 11; Name = "CJ McCollum";               Pos = "PG,SG";      Team = "New Orleans Pelicans" }
    @{ Row = 12; Name = "Zach LaVine";               Pos = "SG,SF";      Team = "Chicago Bulls" }
    @{ Row = 13; Name = "Nikola Jovic";              Pos = "PF,C";       Team = "Miami Heat" }
    @{ Row = 14; Name = "Lauri Markkanen";           Pos = "SF,PF";      Team = "Utah Jazz" }
    @{ Row = 15; Name = "RJ Barrett";                Pos = "SG,SF,PF";   Team = "Toronto Raptors" }
    @{ Row = 16; Name = "Kyrie Irving";              Pos = "PG,SG";      Team = "Dallas Mavericks" }
    @{ Row = 17; Name = "John Collins";              Pos = "PF,C";       Team = "Utah Jazz" }
    @{ Row = 18; Name = "Joel Embiid";                Pos = "C";          Team = "Philadelphia 76ers" }
    @{ Row = 19; Name = "Jimmy Butler";              Pos = "SF,PF";      Team = "Miami Heat" }
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Name
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Pos
    $ws.Cells.Item($entry.Row, 3).Value = $entry.Team
}
